$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Copy() | Out-Null
$ws.Range("J1:J6").PasteSpecial(-4122) | Out-Null

$ws.Range("J1").Value = "GAGATCGTCTCAGGT "
$ws.Range("J1").Font.Size = 12

$ws.Range("J2").Value = "EEEEEEEEEEEEEii"
$ws.Range("J2").Font.Size = 12

$ws.Range("J3").Value = "CAAGATCGGCCCGGT"
$ws.Range("J3").Font.Size = 12

$ws.Range("J4").Value = "iiiiEEEEEEEEEii"
$ws.Range("J4").Font.Size = 12

$ws.Range("J5").Value = "EEEEEEiiiiiEEii"
$ws.Range("J5").Font.Size = 12

$ws.Range("J6").Value = "EEEEiiiiiiiiiii"
$ws.Range("J6").Font.Size = 12

$ws.Columns("J:J").AutoFit() | Out-Null
Write-Host ("Col J width after autofit: " + $ws.Columns("J:J").ColumnWidth())
